$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current")

# --- Rename / re-describe the Hall effect sensor, add part number ---
$ws.Range("A3").Value = "Hall effect sensor"
$ws.Range("B3").Value = "ACS723LLCTR-40AU-T"

# --- New potential divider section (rows 16-19) that feeds the sensor gain ---
$ws.Range("A16").Value = "potential divider"

$ws.Range("A17").Value = "R28"
$ws.Range("B17").Value = 10000
$ws.Range("C17").Value = "ohms"

$ws.Range("A18").Value = "R29"
$ws.Range("B18").Value = 5100
$ws.Range("C18").Value = "ohms"

$ws.Range("A19").Value = "divider factor"
$ws.Range("B19").Formula = "=B17/(B17+B18)"

# --- note about the bias that needs nulling out ---
$ws.Range("A14").Value = "Also note these is a bias to null out."

# --- sensor circuit gain is now driven by the divider factor instead of a fixed 1 ---
$ws.Range("B9").Formula = "=B19"

# --- fix the selection on the sheet view (was A15:B15) ---
$ws.Range("A15").Select()

$wb.Save()
